# Updated/re-processed pipeline output for "signal segment 4" (row 5) across
# the workbook's sheets: the per-sample intensity values on Step1_Data were
# recomputed (mounted pipeline rerun), which cascades into the cumulative
# values on Step2_Sj and the derived threshold statistics on the
# Step3_DataPts_* sheets.
$wb = $excel.ActiveWorkbook

# --- Step1_Data: recomputed per-sample intensity values for row 5 ---
$ws = $wb.Worksheets.Item("Step1_Data")
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0.1237722762792683
$ws.Range("F5").Value = 0.2194619647783052
$ws.Range("G5").Value = 0.09451484611487088
$ws.Range("H5").Value = 0.01320318513425232
$ws.Range("I5").Value = 0.00809148154094954
$ws.Range("M5").Value = 0.03640292962092147
$ws.Range("N5").Value = 0.006772103757577015
$ws.Range("O5").Value = 0.09825555113042249
$ws.Range("P5").Value = 0.01024308869392187
$ws.Range("R5").Value = 0.01451361841226421
$ws.Range("S5").Value = 0.02848539806708211
$ws.Range("T5").Value = 0.06726514825316722
$ws.Range("U5").Value = 0.01972394206323638
$ws.Range("V5").Value = 0.001824588948671209
$ws.Range("W5").Value = 0.0502953496639682
$ws.Range("X5").Value = 0.007004198872256028
$ws.Range("Y5").Value = 0.02547883591338231
$ws.Range("Z5").Value = 0.04872080211224116
$ws.Range("AA5").Value = 0.05233435402739666
$ws.Range("AB5").Value = 0.002414487520262807
$ws.Range("AC5").Value = 0.01856913200669596
$ws.Range("AD5").Value = 0.0161169506933235
$ws.Range("AE5").Value = 0.0006171732326553983
$ws.Range("AF5").Value = 0.022636433779896
$ws.Range("AG5").Value = 0.0009638875318230338
$ws.Range("AI5").Value = 0.007385516021876887
$ws.Range("AJ5").Value = 0.004932755829311757

# --- Step2_Sj: recomputed running cumulative sum of Step1_Data row 5 ---
$ws = $wb.Worksheets.Item("Step2_Sj")
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0.1237722762792683
$ws.Range("F5").Value = 0.3432342410575736
$ws.Range("G5").Value = 0.4377490871724444
$ws.Range("H5").Value = 0.4509522723066968
$ws.Range("I5").Value = 0.4590437538476463
$ws.Range("J5").Value = 0.4590437538476463
$ws.Range("K5").Value = 0.4590437538476463
$ws.Range("L5").Value = 0.4590437538476463
$ws.Range("M5").Value = 0.4954466834685677
$ws.Range("N5").Value = 0.5022187872261448
$ws.Range("O5").Value = 0.6004743383565673
$ws.Range("P5").Value = 0.6107174270504891
$ws.Range("Q5").Value = 0.6107174270504891
$ws.Range("R5").Value = 0.6252310454627533
$ws.Range("S5").Value = 0.6537164435298355
$ws.Range("T5").Value = 0.7209815917830027
$ws.Range("U5").Value = 0.7407055338462391
$ws.Range("V5").Value = 0.7425301227949103
$ws.Range("W5").Value = 0.7928254724588785
$ws.Range("X5").Value = 0.7998296713311345
$ws.Range("Y5").Value = 0.8253085072445168
$ws.Range("Z5").Value = 0.8740293093567579
$ws.Range("AA5").Value = 0.9263636633841545
$ws.Range("AB5").Value = 0.9287781509044174
$ws.Range("AC5").Value = 0.9473472829111134
$ws.Range("AD5").Value = 0.9634642336044369
$ws.Range("AE5").Value = 0.9640814068370923
$ws.Range("AF5").Value = 0.9867178406169883
$ws.Range("AG5").Value = 0.9876817281488114
$ws.Range("AH5").Value = 0.9876817281488114
$ws.Range("AI5").Value = 0.9950672441706883
$ws.Range("AJ5").Value = 1

# --- Step3_DataPts_0.5: updated threshold-crossing stats for row 5 ---
$ws = $wb.Worksheets.Item("Step3_DataPts_0.5")
$ws.Range("D5").Value = 13
$ws.Range("F5").Value = 0.5022187872261448
$ws.Range("G5").Value = 11

# --- Step3_DataPts_0.7: updated cumulative value at crossing for row 5 ---
$ws = $wb.Worksheets.Item("Step3_DataPts_0.7")
$ws.Range("F5").Value = 0.7209815917830027

# --- Step3_DataPts_0.8: updated threshold-crossing stats for row 5 ---
$ws = $wb.Worksheets.Item("Step3_DataPts_0.8")
$ws.Range("D5").Value = 24
$ws.Range("F5").Value = 0.8253085072445168
$ws.Range("G5").Value = 22

# --- Step3_DataPts_0.9: updated cumulative value at crossing for row 5 ---
$ws = $wb.Worksheets.Item("Step3_DataPts_0.9")
$ws.Range("F5").Value = 0.9263636633841545
